# Add 2022-Q4 data:
#  1. Insert a new "2022-Q4" worksheet (holdings detail) right before the
#     existing "2022-Q3" sheet.
#  2. Insert a new row at the top of the "总计" (summary) sheet's data for
#     the 2022-Q4 quarter, shifting the older quarters down and renumbering
#     the running index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q4" worksheet with the per-fund holdings detail.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("2022-Q3"))
$newSheet.Name = "2022-Q4"

# Re-fetch "2022-Q3" by name now that sheets have shifted position, so the
# copy below pulls from the intended sheet rather than the freshly-added one.
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# Pull over the header row + "A" index-column formatting from the sheet
# this data was modeled after so the new sheet matches the look of its
# siblings (bold / bordered / centered header, bordered index column).
$q3Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$q3Sheet.Range("A2:A3").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'005176"
$newSheet.Range("C2").Value = "富国精准医疗灵活配置混合"
$newSheet.Range("D2").Value = "'35.53"
$newSheet.Range("E2").Value = "'93.98"
$newSheet.Range("F2").Value = "'3.54"
$newSheet.Range("G2").Value = "'1.2578"
$newSheet.Range("H2").Value = 9

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'001900"
$newSheet.Range("C3").Value = "诺安精选价值混合"
$newSheet.Range("D3").Value = "'0.13"
$newSheet.Range("E3").Value = "'62.72"
$newSheet.Range("F3").Value = "'1.70"
$newSheet.Range("G3").Value = "'0.0022"
$newSheet.Range("H3").Value = 10

# ---------------------------------------------------------------------
# 2. Update "总计" summary sheet with the new 2022-Q4 row.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()

# The freshly-inserted row inherits blank/auto formatting on B:D; reset it
# back to "no explicit style" like its siblings, then recreate the bordered
# index-column look by copying it from the (now-shifted) row below.
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 1.26

# Renumber the running index (column A) for the quarters that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
